$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp refresh ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 16:35"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1388283
$ws.Range("C4").Value = 2449
$ws.Range("D4").Value = 262326
$ws.Range("E4").Value = 1043939
$ws.Range("F4").Value = 16501
$ws.Range("G4").Value = 223
$ws.Range("H4").Value = 82018

# --- Alemania (row 10) ---
$ws.Range("B10").Value = 172723
$ws.Range("C10").Value = 147
$ws.Range("E10").Value = 17856
$ws.Range("G10").Value = 6
$ws.Range("H10").Value = 7667

# --- India (row 15) ---
$ws.Range("B15").Value = 72593
$ws.Range("C15").Value = 1825
$ws.Range("D15").Value = 23203
$ws.Range("E15").Value = 47059
$ws.Range("G15").Value = 37
$ws.Range("H15").Value = 2331

# --- Albania / Sri Lanka swap places in the ranking (rows 104-105) ---
# Row 104 now holds Sri Lanka's (updated) figures.
$ws.Range("A104").Value = "Sri Lanka"
$ws.Range("B104").Value = 879
$ws.Range("C104").Value = 16
$ws.Range("D104").Value = 366
$ws.Range("E104").Value = 504
$ws.Range("H104").Value = 9

# Row 105 now holds Albania's figures (its previous, unchanged totals).
$ws.Range("A105").Value = "Albania"
$ws.Range("B105").Value = 876
$ws.Range("C105").Value = 4
$ws.Range("D105").Value = 682
$ws.Range("E105").Value = 163
$ws.Range("H105").Value = 31
